# Update "想去人数" (want-to-go count) figures in column F across the
# three sheets that contain event data: 展览, 演出, and 全部类型.
# 本地生活 has no data rows and is left untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 9139
$ws1.Range("F7").Value  = 1378
$ws1.Range("F8").Value  = 206
$ws1.Range("F9").Value  = 62
$ws1.Range("F10").Value = 99
$ws1.Range("F11").Value = 5936
$ws1.Range("F14").Value = 104
$ws1.Range("F15").Value = 4532
$ws1.Range("F16").Value = 16
$ws1.Range("F17").Value = 166
$ws1.Range("F19").Value = 31
$ws1.Range("F23").Value = 262
$ws1.Range("F25").Value = 2792
$ws1.Range("F26").Value = 129

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 39
$ws2.Range("F3").Value = 41

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 9139
$ws4.Range("F5").Value  = 39
$ws4.Range("F8").Value  = 1378
$ws4.Range("F9").Value  = 206
$ws4.Range("F10").Value = 62
$ws4.Range("F11").Value = 99
$ws4.Range("F12").Value = 5936
$ws4.Range("F15").Value = 104
$ws4.Range("F16").Value = 4532
$ws4.Range("F17").Value = 16
$ws4.Range("F18").Value = 166
$ws4.Range("F20").Value = 31
$ws4.Range("F24").Value = 262
$ws4.Range("F26").Value = 2792
$ws4.Range("F27").Value = 41
$ws4.Range("F28").Value = 129
